$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# New ballot row for Jay Cohen
$ws.Range("A33").Value = "Jay Cohen"
$ws.Range("C33").Value = "x"
$ws.Range("D33").Value = "x"
$ws.Range("E33").Value = "x"
$ws.Range("I33").Value = "x"
$ws.Range("J33").Value = "x"
$ws.Range("K33").Value = "x"
$ws.Range("O33").Value = "x"
$ws.Range("P33").Value = "x"
$ws.Range("Q33").Value = "x"
$ws.Range("V33").Value = "x"
$ws.Range("AK33").Value = 10
$ws.Range("AL33").Value = "not specified"
$ws.Range("AM33").Value = 43444
$ws.Range("AM33").NumberFormat = "m/d/yy"

# Update the active selection to the newly added row, matching the
# author's final cursor position after entering the data.
$ws.Activate()
$ws.Range("AL33").Select()
